{"js": "// Insert a new paragraph (\"Hi this is nischl.\") right after the paragraph\n// that contains the ellipsis (\"\u2026\") and before the trailing empty paragraph,\n// reproducing the exact run/proofErr structure Word generates when it\n// flags \"nischl\" as a misspelling.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"\\u2026\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  // Fall back to the last non-empty paragraph (just before the trailing\n  // blank paragraph) if the ellipsis text could not be matched exactly.\n  for (let i = paragraphs.items.length - 1; i >= 0; i--) {\n    if (paragraphs.items[i].text.trim() !== \"\") {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the ellipsis paragraph\");\n}\n\nconst newParagraphOoxml =\n  '<w:p><w:r><w:t xml:space=\"preserve\">Hi this is </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>nischl</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>.</w:t></w:r></w:p>';\n\nconst packagedOoxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newParagraphOoxml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// A zero-length range positioned immediately after the ellipsis paragraph's\n// mark; inserting \"after\" it creates a brand-new sibling paragraph rather\n// than merging content into the existing one.\nconst insertionPoint = target.getRange(Word.RangeLocation.after);\ninsertionPoint.insertOoxml(packagedOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new paragraph (\"Hi this is nischl.\") right after the paragraph\n# that holds the ellipsis (\"\u2026\") and before the final, empty paragraph that\n# precedes the section break.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is just the ellipsis character + mark.\n$ellipsis = [char]0x2026\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text.Trim() -eq $ellipsis) {\n        $target = $para\n        break\n    }\n}\nif (-not $target) {\n    # Fall back to the last non-empty paragraph (just before the trailing\n    # blank paragraph) if the ellipsis text could not be matched exactly.\n    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n        $para = $d.Paragraphs.Item($i)\n        if ($para.Range.Text.Trim() -ne \"\") {\n            $target = $para\n            break\n        }\n    }\n}\nif (-not $target) {\n    throw \"Could not locate the ellipsis paragraph\"\n}\n\n# Create a brand-new empty paragraph right after it.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n\n# Populate the new paragraph with the exact run/proofErr structure using\n# WordprocessingML so the misspelling marker around \"nischl\" is preserved.\n$newParaXml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">Hi this is </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>nischl</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p>\n'@\n\n$newPara.Range.InsertXML($newParaXml)\n"}
